$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update existing row 59 (2025-07): retained_customers and retention_rate changed
$ws.Range("B59").Value = 157
$ws.Range("D59").Value = 68.5589519650655

# Add new row 60 (2025-08)
$ws.Range("A60").Value = "2025-08"
$ws.Range("B60").Value = 11
$ws.Range("C60").Value = 238
$ws.Range("D60").Value = 4.621848739495799
